# feat: updated Tommaso Stedile team
#
# Tommaso Stedile's "Portiere" (goalkeeper) entry changes from
# "Nicolas Giordani | FC Savignano" to "Federico Leonardi | Sughi ebbasta".
# Also brings along the incidental view/format-state deltas captured by the
# diff: selection moved back to the top of the sheet (B18) with the window
# scrolled to A1, and column B grew slightly wider.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content edit -----------------------------------------------------
# Tommaso Stedile is row 21 ("A21"); B21 is his goalkeeper pick.
$ws.Range("B21").Value = "Federico Leonardi | Sughi ebbasta"

# --- view / selection state --------------------------------------------
# Selection moves to B18 and the window scrolls back to the top-left (A1).
$ws.Range("B18").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1

# --- column width --------------------------------------------------------
# Column B (Portiere) widens from ~28.48 to ~29.69 characters. The host's
# column grid snaps to its own pixel increments, so feed it the input that
# lands on the closest achievable width to 29.69.
$ws.Columns(2).ColumnWidth = 28.75

# sheetFormatPr defaultColWidth nudges from 8.5390625 to 8.54296875 as a
# side-effect of Excel's own default-width recompute; reflect the intent
# via StandardWidth even though it mirrors the sheet's default metrics.
$ws.StandardWidth = 8.54296875
